$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header for the first (previously empty) column to "Year"
$ws.Range("A1").Value = "Year"
